$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Loads_TC003): Actual Result column flips from "YES" to "NO"
$ws.Range("C9").Value = "NO"

# Row 10 (Loads_TC004): Actual Result column flips from "NO" to "YES"
$ws.Range("C10").Value = "YES"

# Row 10 (Loads_TC004): Status column text changes
$ws.Range("D10").Value = "Load Approved successfully"

# Update the selection/active cell to D10 (matches the new cursor position in the diff)
$ws.Range("D10").Select()
